$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-05 Saturday" "2025-04-06 Sunday"

Replace-Text "79×73=" "89×58="
Replace-Text "17×70=" "84×41="
Replace-Text "92×17=" "54×55="
Replace-Text "32×86=" "18×77="
Replace-Text "73×65=" "31×21="

Replace-Text "16×89=" "90×42="
Replace-Text "52×71=" "83×59="
Replace-Text "31×92=" "15×49="
Replace-Text "70×38=" "25×90="
Replace-Text "97×92=" "24×46="

Replace-Text "82×48=" "83×25="
Replace-Text "42×50=" "87×31="
Replace-Text "11×62=" "11×54="
Replace-Text "13×32=" "47×95="
Replace-Text "78×76=" "16×99="

Replace-Text "84×54=" "59×35="
Replace-Text "25×42=" "24×93="
Replace-Text "79×49=" "17×35="
Replace-Text "21×94=" "43×67="
Replace-Text "74×20=" "33×58="

Replace-Text "14×98=" "40×22="
Replace-Text "65×49=" "46×20="
Replace-Text "17×27=" "58×98="
Replace-Text "92×43=" "51×80="
Replace-Text "94×51=" "22×51="
